$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text so numeric-looking strings
# (e.g. "1.011") are not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '27.303.53'
$ws.Range('E2').Value = '  +1.30%  '
$ws.Range('D3').Value = '1.833.16'
$ws.Range('E3').Value = '  +0.82%  '
$ws.Range('D4').Value = '1.011'
$ws.Range('E4').Value = '  +0.85%  '
$ws.Range('D5').Value = '314.69'
$ws.Range('E5').Value = '  +1.62%  '
$ws.Range('E6').Value = '  +0.79%  '
$ws.Range('E7').Value = '  +1.77%  '
$ws.Range('D8').Value = '0.3688'
$ws.Range('E8').Value = '  +0.69%  '
$ws.Range('E9').Value = '  +1.00%  '
$ws.Range('D10').Value = '0.8862'
$ws.Range('E10').Value = '  +1.66%  '
$ws.Range('E11').Value = '  +1.08%  '
$ws.Range('D12').Value = '1.884.20'
$ws.Range('E12').Value = '  +3.73%  '
$ws.Range('D13').Value = '0.07333'
$ws.Range('E13').Value = '  +3.08%  '
$ws.Range('D14').Value = '5.433'
$ws.Range('E14').Value = '  +0.83%  '
$ws.Range('D15').Value = '94.03'
$ws.Range('E15').Value = '  +2.97%  '
$ws.Range('D16').Value = '6.563'
$ws.Range('E16').Value = '  +0.77%  '
$ws.Range('E17').Value = '  +0.59%  '
$ws.Range('D18').Value = '0.000008788'
$ws.Range('E18').Value = '  +1.05%  '
$ws.Range('E19').Value = '  +0.87%  '
$ws.Range('D20').Value = '27.539.21'
$ws.Range('E20').Value = '  +2.10%  '
$ws.Range('D21').Value = '14.79'
$ws.Range('E21').Value = '  +0.92%  '
$ws.Range('D22').Value = '5.287'
$ws.Range('D23').Value = '10.67'
$ws.Range('E23').Value = '  +0.80%  '
$ws.Range('D24').Value = '2.096.87'
$ws.Range('E24').Value = '  +2.67%  '
$ws.Range('D25').Value = '1.892'
$ws.Range('E25').Value = '  -0.15%  '
$ws.Range('E26').Value = '  +0.55%  '
$ws.Range('E27').Value = '  +1.18%  '
$ws.Range('D28').Value = '2.144'
$ws.Range('E28').Value = '  +0.18%  '
$ws.Range('D29').Value = '5.235'
$ws.Range('E29').Value = '  -0.53%  '
$ws.Range('D30').Value = '117.24'
$ws.Range('E30').Value = '  +0.69%  '
$ws.Range('D31').Value = '0.08994'
$ws.Range('E31').Value = '  +1.12%  '
$ws.Range('D32').Value = '0.7493'
$ws.Range('E32').Value = '  -1.37%  '
$ws.Range('D33').Value = '1.175'
$ws.Range('E33').Value = '  +0.74%  '
$ws.Range('D34').Value = '4.546'
$ws.Range('E34').Value = '  +1.29%  '
$ws.Range('E35').Value = '  +1.39%  '
$ws.Range('E36').Value = '  +0.92%  '
$ws.Range('E37').Value = '  -0.08%  '
$ws.Range('D38').Value = '0.05347'
$ws.Range('E38').Value = '  +1.12%  '
$ws.Range('D39').Value = '0.01955'
$ws.Range('E39').Value = '  +0.44%  '
$ws.Range('D40').Value = '2.977'
$ws.Range('E40').Value = '  -0.17%  '
$ws.Range('D41').Value = '2.395'
$ws.Range('E41').Value = '  +2.85%  '
$ws.Range('D42').Value = '7.243'
$ws.Range('E42').Value = '  +1.31%  '
$ws.Range('E43').Value = '  -0.09%  '
$ws.Range('D44').Value = '0.1658'
$ws.Range('E44').Value = '  +0.11%  '
$ws.Range('D45').Value = '8.482'
$ws.Range('E45').Value = '  +0.59%  '
$ws.Range('E46').Value = '  +1.61%  '
$ws.Range('E47').Value = '  +0.93%  '
$ws.Range('E48').Value = '  +0.88%  '
$ws.Range('D49').Value = '104.99'
$ws.Range('E49').Value = '  +1.55%  '
$ws.Range('D50').Value = '1.669'
$ws.Range('E50').Value = '  +0.42%  '
$ws.Range('D51').Value = '0.06296'
$ws.Range('E51').Value = '  +0.08%  '

# Restore default (unformatted) styling so the text-format override
# doesn't linger as a visible style change on the cells.
$ws.Range("D2:D51").ClearFormats()
